# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 88 (pushing existing rows 88-106 down to 89-107)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 88, shifting data down
$ws.Rows.Item(88).Insert()

# Populate the newly inserted row 88 with the new weekly record
$ws.Range("A88").Value = 10
$ws.Range("B88").Value = "Vega Modelo de Temuco"
$ws.Range("C88").Value = "La Araucanía"
$ws.Range("D88").Value = 44637
$ws.Range("E88").Value = 9
$ws.Range("F88").Value = "Fruta"
$ws.Range("G88").Value = 100104
$ws.Range("H88").Value = "Frutos de pepita"
$ws.Range("I88").Value = 100104001
$ws.Range("J88").Value = "Granada"
$ws.Range("K88").Value = "Sin especificar"
$ws.Range("L88").Value = "Primera"
$ws.Range("M88").Value = 25
$ws.Range("N88").Value = 20000
$ws.Range("O88").Value = 20000
$ws.Range("P88").Value = 20000
$ws.Range("Q88").Value = "$/caja 15 kilos empedrada"
$ws.Range("R88").Value = "Provincia de Limarí"
$ws.Range("S88").Value = 1333
$ws.Range("T88").Value = 15

# Ensure the date cell keeps the same date number format as the other date cells in column D
$ws.Range("D88").NumberFormat = $ws.Range("D89").NumberFormat
